$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Source: Deshmuk 2015" note from row 11 up to row 9 ---
$ws.Range("A9").Value = "Source: Deshmuk 2015"
$ws.Range("A11").ClearContents()

# --- New attribution line for the ELAVI-provided data (added to the
#     shared-string table before "Utilities" so the indices line up) ---
$ws.Range("A19").Value = "Source: ELAVI"

# --- New column D (Utilities) for the first table (rows 1-6) ---
$ws.Range("D1").Value = "Utilities"
$ws.Range("D2").Value = 0.94
$ws.Range("D3").Value = 0.87
$ws.Range("D4").Value = 0.8
$ws.Range("D6").Formula = "=B2*D2+B3*D3+B4*D4"
$ws.Range("D6").Interior.Color = $ws.Range("C6").Interior.Color

# Highlight the "Normal to HGAIN" probability column (B) with the new
# light "Background 2" themed fill
$ws.Range("B2:B4").Interior.ThemeColor = 4

# --- New separator row 10 (solid black fill across A:C) ---
$ws.Range("A10:C10").Interior.ThemeColor = 1

# --- New second table (rows 12-17), the ELAVI-calibrated version ---
$ws.Range("B12").Value = "CD4 distribution"
$ws.Range("C12").Value = "Normal to HGAIN"
$ws.Range("D12").Value = "Utilities"

$ws.Range("A13").Value = ">500"
$ws.Range("B13").Value = 0.83898305084745761
$ws.Range("C13").Value = 0.053
$ws.Range("D13").Value = 0.94

$ws.Range("A14").Value = "200-500"
$ws.Range("B14").Value = 0.14689265536723164
$ws.Range("C14").Value = 0.053
$ws.Range("D14").Value = 0.87

$ws.Range("A15").Value = "<200"
$ws.Range("B15").Value = 0.014124293785310734
$ws.Range("C15").Value = 0.138
$ws.Range("D15").Value = 0.8

$ws.Range("B13:B15").Interior.ThemeColor = 4

$ws.Range("C17").Formula = "=B13*C13+B14*C14+B15*C15"
$ws.Range("D17").Formula = "=B13*D13+B14*D14+B15*D15"
$ws.Range("C17:D17").Interior.Color = $ws.Range("C6").Interior.Color

# --- Selection / view tidy-up to match what was recorded ---
$ws.Range("C9").Select()
